# updated EEG Subsidy: Found a way to properly isolate the subsidy payment.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "estimate" row (row 2) now reports the isolated subsidy payment per MWh
# of CO2 offset, rather than being labelled as a generic "denominator".
$ws.Range("A2").Value = "co2_offset_per_MWh"
$ws.Range("B2").Value = -233.2
$ws.Range("C2").Value = 10.1

# Set print/page layout to A4 portrait.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Move the active cell selection.
$ws.Range("F5").Select() | Out-Null
